$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-label / reword the note column and text updates -------------------
# Move "note" header from E4 to E3 (leaving E4 blank)
$ws.Range("E4").Value = ""
$ws.Range("E3").Value = "note"

# Reword the notes in column E (rows 5-7) to reflect the new logistic-function
# based step input and the resulting angular velocity figures.
$ws.Range("E5").Value = "angular velocity 8000 during response to 0.10"
$ws.Range("E6").Value = "PID Tuner for 2 response time and 0.6 robustness, but the issue is its tuning for a step t o 1, and I'm only stepping to 0.1"
$ws.Range("E7").Value = "based on advice from James Viollete, 1000 angular velocity during step respons"

# Rename the "data" header (A4) to "date"
$ws.Range("A4").Value = "date"

# --- Swap the recorded p/i/d values between rows 6 and 7 -------------------
$ws.Range("B6").Value = -0.032627999999999997
$ws.Range("C6").Value = -0.0018881
$ws.Range("D6").Value = -0.14096

$ws.Range("B7").Value = -0.01
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -1

# --- Bold the header rows ---------------------------------------------------
$ws.Range("B3").Font.Bold = $true
$ws.Range("E3").Font.Bold = $true
$ws.Range("A4:D4").Font.Bold = $true

# --- Page setup (portrait orientation) -------------------------------------
$ws.PageSetup.Orientation = 1
